$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 265.0625
$ws.Range("I33").Value = 242.6
$ws.Range("K33").Value = 242.6
$ws.Range("M33").Value = -13.59999999999999
$ws.Range("H112").Value = 1275.8689
$ws.Range("J112").Value = 1297.0847
$ws.Range("L112").Value = 3891.2541
$ws.Range("N112").Value = -6107.2541
$ws.Range("H113").Value = 3669.6875
$ws.Range("I113").Value = 2265.3572
$ws.Range("K113").Value = 2265.3572
$ws.Range("M113").Value = 988.6428000000001
$ws.Range("H133").Value = 55833.332
$ws.Range("J133").Value = 55833.332
$ws.Range("L133").Value = 55833.332
$ws.Range("N133").Value = -65953.33199999999
$ws.Range("H137").Value = 1254894
$ws.Range("I137").Value = 1537482.9
$ws.Range("J137").Value = 3428.5715
$ws.Range("K137").Value = 4612448.699999999
$ws.Range("L137").Value = 10285.7145
$ws.Range("M137").Value = -4609898.699999999
$ws.Range("N137").Value = -15385.7145
$ws.Range("H138").Value = 2930.4546
$ws.Range("I138").Value = 2088
$ws.Range("J138").Value = 3070.8635
$ws.Range("K138").Value = 6264
$ws.Range("L138").Value = 9212.5905
$ws.Range("M138").Value = -1124
$ws.Range("N138").Value = -19492.5905
$ws.Range("H141").Value = 156090.92
$ws.Range("I141").Value = 201608.7
$ws.Range("K141").Value = 604826.1000000001
$ws.Range("M141").Value = -599646.1000000001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 816.6667
$ws.Range("I2").Value = 700
$ws.Range("J2").Value = 1050
$ws.Range("K2").Value = 700
$ws.Range("L2").Value = 1050
$ws.Range("M2").Value = -587
$ws.Range("N2").Value = -1276
$ws.Range("H45").Value = 1196.7
$ws.Range("I45").Value = 1245.875
$ws.Range("J45").Value = 1000
$ws.Range("K45").Value = 1245.875
$ws.Range("L45").Value = 1000
$ws.Range("M45").Value = -868.875
$ws.Range("N45").Value = -1754
$ws.Range("H116").Value = 816.6667
$ws.Range("I116").Value = 700
$ws.Range("J116").Value = 1050
$ws.Range("K116").Value = 700
$ws.Range("L116").Value = 1050
$ws.Range("M116").Value = 1594
$ws.Range("N116").Value = -5638

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 816.6667
$ws.Range("I3").Value = 700
$ws.Range("J3").Value = 1050
$ws.Range("K3").Value = 700
$ws.Range("L3").Value = 1050
$ws.Range("M3").Value = -586
$ws.Range("N3").Value = -1278
$ws.Range("H59").Value = 45260
$ws.Range("J59").Value = 45260
$ws.Range("L59").Value = 45260
$ws.Range("N59").Value = -46954
$ws.Range("H137").Value = 45085.715
$ws.Range("J137").Value = 45085.715
$ws.Range("L137").Value = 45085.715
$ws.Range("N137").Value = -55285.715
$ws.Range("H138").Value = 40890.742
$ws.Range("J138").Value = 40890.742
$ws.Range("L138").Value = 40890.742
$ws.Range("N138").Value = -51170.742
$ws.Range("H140").Value = 50528
$ws.Range("J140").Value = 50528
$ws.Range("L140").Value = 50528
$ws.Range("N140").Value = -60888

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 10101538
$ws.Range("I16").Value = 18518982
$ws.Range("J16").Value = 605.2
$ws.Range("K16").Value = 18518982
$ws.Range("L16").Value = 605.2
$ws.Range("M16").Value = -18518695
$ws.Range("N16").Value = -1179.2
$ws.Range("H31").Value = 3204.5557
$ws.Range("I31").Value = 1356.2778
$ws.Range("J31").Value = 6901.1113
$ws.Range("K31").Value = 1356.2778
$ws.Range("L31").Value = 6901.1113
$ws.Range("M31").Value = -1061.2778
$ws.Range("N31").Value = -7491.1113
$ws.Range("H34").Value = 3204.5557
$ws.Range("I34").Value = 1356.2778
$ws.Range("J34").Value = 6901.1113
$ws.Range("K34").Value = 1356.2778
$ws.Range("L34").Value = 6901.1113
$ws.Range("M34").Value = -1154.2778
$ws.Range("N34").Value = -7305.1113
$ws.Range("H59").Value = 37625.5
$ws.Range("J59").Value = 37625.5
$ws.Range("L59").Value = 37625.5
$ws.Range("N59").Value = -39915.5
$ws.Range("H74").Value = 33701
$ws.Range("J74").Value = 33701
$ws.Range("L74").Value = 33701
$ws.Range("N74").Value = -35449
$ws.Range("H77").Value = 33701
$ws.Range("J77").Value = 33701
$ws.Range("L77").Value = 101103
$ws.Range("N77").Value = -109839
$ws.Range("H113").Value = 10101538
$ws.Range("I113").Value = 18518982
$ws.Range("J113").Value = 605.2
$ws.Range("K113").Value = 18518982
$ws.Range("L113").Value = 605.2
$ws.Range("M113").Value = -18516812
$ws.Range("N113").Value = -4945.2
$ws.Range("H137").Value = 32751.25
$ws.Range("J137").Value = 32751.25
$ws.Range("L137").Value = 32751.25
$ws.Range("N137").Value = -42951.25
$ws.Range("H138").Value = 44500
$ws.Range("J138").Value = 44500
$ws.Range("L138").Value = 44500
$ws.Range("N138").Value = -54780
$ws.Range("H140").Value = 113672.5
$ws.Range("J140").Value = 113672.5
$ws.Range("L140").Value = 113672.5
$ws.Range("N140").Value = -124032.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 1239.6923
$ws.Range("I109").Value = 1123.8334
$ws.Range("J109").Value = 2630
$ws.Range("K109").Value = 3371.5002
$ws.Range("L109").Value = 7890
$ws.Range("M109").Value = -2331.5002
$ws.Range("N109").Value = -9970
$ws.Range("H112").Value = 6060
$ws.Range("I112").Value = 3433.3333
$ws.Range("K112").Value = 10299.9999
$ws.Range("M112").Value = -9191.999899999999
$ws.Range("H113").Value = 3378961
$ws.Range("I113").Value = 586.0741
$ws.Range("J113").Value = 12500573
$ws.Range("K113").Value = 1758.2223
$ws.Range("L113").Value = 37501719
$ws.Range("M113").Value = 411.7776999999999
$ws.Range("N113").Value = -37506059
$ws.Range("H116").Value = 804.75
$ws.Range("I116").Value = 804.75
$ws.Range("K116").Value = 2414.25
$ws.Range("M116").Value = 1027.75
$ws.Range("H117").Value = 2375
$ws.Range("I117").Value = 1500
$ws.Range("J117").Value = 2666.6667
$ws.Range("K117").Value = 4500
$ws.Range("L117").Value = 8000.000100000001
$ws.Range("M117").Value = -1058
$ws.Range("N117").Value = -14884.0001
$ws.Range("H118").Value = 3708.8333
$ws.Range("I118").Value = 563.25
$ws.Range("K118").Value = 1689.75
$ws.Range("M118").Value = -446.75
$ws.Range("H121").Value = 1980.125
$ws.Range("I121").Value = 349.875
$ws.Range("J121").Value = 2251.8333
$ws.Range("K121").Value = 1049.625
$ws.Range("L121").Value = 6755.499899999999
$ws.Range("M121").Value = 260.375
$ws.Range("N121").Value = -9375.499899999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2674.7273
$ws.Range("I113").Value = 2990.25
$ws.Range("J113").Value = 1833.3334
$ws.Range("K113").Value = 2990.25
$ws.Range("L113").Value = 1833.3334
$ws.Range("M113").Value = -820.25
$ws.Range("N113").Value = -6173.3334
$ws.Range("H137").Value = 40500
$ws.Range("J137").Value = 40500
$ws.Range("L137").Value = 40500
$ws.Range("N137").Value = -50700
$ws.Range("H140").Value = 38309.547
$ws.Range("J140").Value = 38309.547
$ws.Range("L140").Value = 38309.547
$ws.Range("N140").Value = -48669.547

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 17885.834
$ws.Range("J43").Value = 17885.834
$ws.Range("L43").Value = 17885.834
$ws.Range("N43").Value = -18271.834
$ws.Range("H64").Value = 37250
$ws.Range("J64").Value = 37250
$ws.Range("L64").Value = 37250
$ws.Range("N64").Value = -37700
$ws.Range("H67").Value = 37250
$ws.Range("J67").Value = 37250
$ws.Range("L67").Value = 37250
$ws.Range("N67").Value = -38810
$ws.Range("H92").Value = 37500
$ws.Range("J92").Value = 37500
$ws.Range("L92").Value = 37500
$ws.Range("N92").Value = -42492
$ws.Range("H122").Value = 8448.833000000001
$ws.Range("I122").Value = 7296.4
$ws.Range("K122").Value = 21889.2
$ws.Range("M122").Value = -19439.2
$ws.Range("H141").Value = 32119.736
$ws.Range("J141").Value = 32119.736
$ws.Range("L141").Value = 32119.736
$ws.Range("N141").Value = -42479.736

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 67050.89999999999
$ws.Range("J46").Value = 67050.89999999999
$ws.Range("L46").Value = 67050.89999999999
$ws.Range("N46").Value = -67512.89999999999
$ws.Range("H63").Value = 28600
$ws.Range("J63").Value = 28600
$ws.Range("L63").Value = 28600
$ws.Range("N63").Value = -29848
$ws.Range("H66").Value = 28600
$ws.Range("J66").Value = 28600
$ws.Range("L66").Value = 85800
$ws.Range("N66").Value = -92040
$ws.Range("H122").Value = 5306.1035
$ws.Range("I122").Value = 3772.5264
$ws.Range("J122").Value = 8219.9
$ws.Range("K122").Value = 11317.5792
$ws.Range("L122").Value = 24659.7
$ws.Range("M122").Value = -8867.5792
$ws.Range("N122").Value = -29559.7
$ws.Range("H134").Value = 67050.89999999999
$ws.Range("J134").Value = 67050.89999999999
$ws.Range("L134").Value = 201152.7
$ws.Range("N134").Value = -206222.7
$ws.Range("H135").Value = 51810
$ws.Range("J135").Value = 51810
$ws.Range("L135").Value = 51810
$ws.Range("N135").Value = -61950
